$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-04-02 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-04-03 Thursday", 2) | Out-Null
$d.Content.Find.Execute("95-77=", $true, $false, $false, $false, $false, $true, 1, $false, "32+48=", 2) | Out-Null
$d.Content.Find.Execute("38-32=", $true, $false, $false, $false, $false, $true, 1, $false, "26-13=", 2) | Out-Null
$d.Content.Find.Execute("33+25=", $true, $false, $false, $false, $false, $true, 1, $false, "69+4=", 2) | Out-Null
$d.Content.Find.Execute("47-26=", $true, $false, $false, $false, $false, $true, 1, $false, "45-1=", 2) | Out-Null
$d.Content.Find.Execute("25+32=", $true, $false, $false, $false, $false, $true, 1, $false, "30+52=", 2) | Out-Null
$d.Content.Find.Execute("63-59=", $true, $false, $false, $false, $false, $true, 1, $false, "98-3=", 2) | Out-Null
$d.Content.Find.Execute("20+26=", $true, $false, $false, $false, $false, $true, 1, $false, "64-41=", 2) | Out-Null
$d.Content.Find.Execute("63-22=", $true, $false, $false, $false, $false, $true, 1, $false, "82-42=", 2) | Out-Null
$d.Content.Find.Execute("19+57=", $true, $false, $false, $false, $false, $true, 1, $false, "26+17=", 2) | Out-Null
$d.Content.Find.Execute("53+29=", $true, $false, $false, $false, $false, $true, 1, $false, "80+3=", 2) | Out-Null
$d.Content.Find.Execute("8+19=", $true, $false, $false, $false, $false, $true, 1, $false, "58+32=", 2) | Out-Null
$d.Content.Find.Execute("67+27=", $true, $false, $false, $false, $false, $true, 1, $false, "46+16=", 2) | Out-Null
$d.Content.Find.Execute("39-19=", $true, $false, $false, $false, $false, $true, 1, $false, "35+58=", 2) | Out-Null
$d.Content.Find.Execute("10+2=", $true, $false, $false, $false, $false, $true, 1, $false, "72+14=", 2) | Out-Null
$d.Content.Find.Execute("89-67=", $true, $false, $false, $false, $false, $true, 1, $false, "15+7=", 2) | Out-Null
$d.Content.Find.Execute("44-41=", $true, $false, $false, $false, $false, $true, 1, $false, "89-68=", 2) | Out-Null
$d.Content.Find.Execute("75-52=", $true, $false, $false, $false, $false, $true, 1, $false, "18+72=", 2) | Out-Null
$d.Content.Find.Execute("66-24=", $true, $false, $false, $false, $false, $true, 1, $false, "25+10=", 2) | Out-Null
$d.Content.Find.Execute("46-33=", $true, $false, $false, $false, $false, $true, 1, $false, "34+4=", 2) | Out-Null
$d.Content.Find.Execute("59-11=", $true, $false, $false, $false, $false, $true, 1, $false, "43-26=", 2) | Out-Null
$d.Content.Find.Execute("47-46=", $true, $false, $false, $false, $false, $true, 1, $false, "16+19=", 2) | Out-Null
$d.Content.Find.Execute("80-34=", $true, $false, $false, $false, $false, $true, 1, $false, "85-84=", 2) | Out-Null
$d.Content.Find.Execute("44+6=", $true, $false, $false, $false, $false, $true, 1, $false, "84-23=", 2) | Out-Null
$d.Content.Find.Execute("3+59=", $true, $false, $false, $false, $false, $true, 1, $false, "18+60=", 2) | Out-Null
$d.Content.Find.Execute("76-72=", $true, $false, $false, $false, $false, $true, 1, $false, "4+62=", 2) | Out-Null
$d.Content.Find.Execute("85-37=", $true, $false, $false, $false, $false, $true, 1, $false, "99-56=", 2) | Out-Null
$d.Content.Find.Execute("33-27=", $true, $false, $false, $false, $false, $true, 1, $false, "8+56=", 2) | Out-Null
$d.Content.Find.Execute("55+3=", $true, $false, $false, $false, $false, $true, 1, $false, "47-42=", 2) | Out-Null
$d.Content.Find.Execute("97-87=", $true, $false, $false, $false, $false, $true, 1, $false, "89-17=", 2) | Out-Null
$d.Content.Find.Execute("13+21=", $true, $false, $false, $false, $false, $true, 1, $false, "2+16=", 2) | Out-Null
$d.Content.Find.Execute("67-61=", $true, $false, $false, $false, $false, $true, 1, $false, "34+26=", 2) | Out-Null
$d.Content.Find.Execute("21-15=", $true, $false, $false, $false, $false, $true, 1, $false, "20+22=", 2) | Out-Null
$d.Content.Find.Execute("98-96=", $true, $false, $false, $false, $false, $true, 1, $false, "26+57=", 2) | Out-Null
$d.Content.Find.Execute("1+87=", $true, $false, $false, $false, $false, $true, 1, $false, "39-14=", 2) | Out-Null
$d.Content.Find.Execute("86-35=", $true, $false, $false, $false, $false, $true, 1, $false, "98-59=", 2) | Out-Null
$d.Content.Find.Execute("43-9=", $true, $false, $false, $false, $false, $true, 1, $false, "78-12=", 2) | Out-Null
$d.Content.Find.Execute("10+32=", $true, $false, $false, $false, $false, $true, 1, $false, "6+66=", 2) | Out-Null
$d.Content.Find.Execute("72-29=", $true, $false, $false, $false, $false, $true, 1, $false, "58+21=", 2) | Out-Null
$d.Content.Find.Execute("55-9=", $true, $false, $false, $false, $false, $true, 1, $false, "83-23=", 2) | Out-Null
$d.Content.Find.Execute("48+14=", $true, $false, $false, $false, $false, $true, 1, $false, "78-56=", 2) | Out-Null
$d.Content.Find.Execute("78-11=", $true, $false, $false, $false, $false, $true, 1, $false, "62-20=", 2) | Out-Null
$d.Content.Find.Execute("46+22=", $true, $false, $false, $false, $false, $true, 1, $false, "80+4=", 2) | Out-Null
$d.Content.Find.Execute("86-57=", $true, $false, $false, $false, $false, $true, 1, $false, "35-0=", 2) | Out-Null
$d.Content.Find.Execute("35-23=", $true, $false, $false, $false, $false, $true, 1, $false, "59-51=", 2) | Out-Null
$d.Content.Find.Execute("13+27=", $true, $false, $false, $false, $false, $true, 1, $false, "92-41=", 2) | Out-Null
$d.Content.Find.Execute("51-48=", $true, $false, $false, $false, $false, $true, 1, $false, "56-4=", 2) | Out-Null
$d.Content.Find.Execute("78+20=", $true, $false, $false, $false, $false, $true, 1, $false, "22+72=", 2) | Out-Null
$d.Content.Find.Execute("78-27=", $true, $false, $false, $false, $false, $true, 1, $false, "59-4=", 2) | Out-Null
$d.Content.Find.Execute("91-46=", $true, $false, $false, $false, $false, $true, 1, $false, "98-86=", 2) | Out-Null
$d.Content.Find.Execute("87-52=", $true, $false, $false, $false, $false, $true, 1, $false, "76-22=", 2) | Out-Null
$d.Content.Find.Execute("79+19=", $true, $false, $false, $false, $false, $true, 1, $false, "94+3=", 2) | Out-Null
$d.Content.Find.Execute("23-0=", $true, $false, $false, $false, $false, $true, 1, $false, "32+38=", 2) | Out-Null
$d.Content.Find.Execute("56-23=", $true, $false, $false, $false, $false, $true, 1, $false, "94-23=", 2) | Out-Null
$d.Content.Find.Execute("57+9=", $true, $false, $false, $false, $false, $true, 1, $false, "74-23=", 2) | Out-Null
$d.Content.Find.Execute("95-83=", $true, $false, $false, $false, $false, $true, 1, $false, "81-20=", 2) | Out-Null
$d.Content.Find.Execute("20-14=", $true, $false, $false, $false, $false, $true, 1, $false, "81-4=", 2) | Out-Null
$d.Content.Find.Execute("53-7=", $true, $false, $false, $false, $false, $true, 1, $false, "57-17=", 2) | Out-Null
$d.Content.Find.Execute("86-43=", $true, $false, $false, $false, $false, $true, 1, $false, "55-18=", 2) | Out-Null
$d.Content.Find.Execute("54-49=", $true, $false, $false, $false, $false, $true, 1, $false, "15+80=", 2) | Out-Null
$d.Content.Find.Execute("55-54=", $true, $false, $false, $false, $false, $true, 1, $false, "28-6=", 2) | Out-Null
$d.Content.Find.Execute("5+8=", $true, $false, $false, $false, $false, $true, 1, $false, "38+47=", 2) | Out-Null
$d.Content.Find.Execute("43-32=", $true, $false, $false, $false, $false, $true, 1, $false, "32+22=", 2) | Out-Null
$d.Content.Find.Execute("74+6=", $true, $false, $false, $false, $false, $true, 1, $false, "23+75=", 2) | Out-Null
$d.Content.Find.Execute("65-56=", $true, $false, $false, $false, $false, $true, 1, $false, "12+16=", 2) | Out-Null
$d.Content.Find.Execute("97-24=", $true, $false, $false, $false, $false, $true, 1, $false, "28+29=", 2) | Out-Null
$d.Content.Find.Execute("17+11=", $true, $false, $false, $false, $false, $true, 1, $false, "65-57=", 2) | Out-Null
$d.Content.Find.Execute("54-29=", $true, $false, $false, $false, $false, $true, 1, $false, "54-10=", 2) | Out-Null
$d.Content.Find.Execute("52+25=", $true, $false, $false, $false, $false, $true, 1, $false, "2+42=", 2) | Out-Null
$d.Content.Find.Execute("9+33=", $true, $false, $false, $false, $false, $true, 1, $false, "37+30=", 2) | Out-Null
$d.Content.Find.Execute("77-21=", $true, $false, $false, $false, $false, $true, 1, $false, "37+48=", 2) | Out-Null
$d.Content.Find.Execute("14+50=", $true, $false, $false, $false, $false, $true, 1, $false, "16+43=", 2) | Out-Null
$d.Content.Find.Execute("58-33=", $true, $false, $false, $false, $false, $true, 1, $false, "95-71=", 2) | Out-Null
$d.Content.Find.Execute("48-11=", $true, $false, $false, $false, $false, $true, 1, $false, "44+14=", 2) | Out-Null
$d.Content.Find.Execute("19+16=", $true, $false, $false, $false, $false, $true, 1, $false, "18-5=", 2) | Out-Null
$d.Content.Find.Execute("92-56=", $true, $false, $false, $false, $false, $true, 1, $false, "27-2=", 2) | Out-Null
$d.Content.Find.Execute("79-73=", $true, $false, $false, $false, $false, $true, 1, $false, "60-38=", 2) | Out-Null
$d.Content.Find.Execute("64-6=", $true, $false, $false, $false, $false, $true, 1, $false, "89-47=", 2) | Out-Null
$d.Content.Find.Execute("8+90=", $true, $false, $false, $false, $false, $true, 1, $false, "22-15=", 2) | Out-Null
$d.Content.Find.Execute("50+36=", $true, $false, $false, $false, $false, $true, 1, $false, "91-0=", 2) | Out-Null
$d.Content.Find.Execute("17+9=", $true, $false, $false, $false, $false, $true, 1, $false, "78-36=", 2) | Out-Null
$d.Content.Find.Execute("38+12=", $true, $false, $false, $false, $false, $true, 1, $false, "61+28=", 2) | Out-Null
$d.Content.Find.Execute("64-27=", $true, $false, $false, $false, $false, $true, 1, $false, "63-29=", 2) | Out-Null
$d.Content.Find.Execute("33+7=", $true, $false, $false, $false, $false, $true, 1, $false, "70-43=", 2) | Out-Null
$d.Content.Find.Execute("55-2=", $true, $false, $false, $false, $false, $true, 1, $false, "79-4=", 2) | Out-Null
$d.Content.Find.Execute("24+26=", $true, $false, $false, $false, $false, $true, 1, $false, "9+57=", 2) | Out-Null
$d.Content.Find.Execute("30+4=", $true, $false, $false, $false, $false, $true, 1, $false, "29+13=", 2) | Out-Null
$d.Content.Find.Execute("49-8=", $true, $false, $false, $false, $false, $true, 1, $false, "80-20=", 2) | Out-Null
$d.Content.Find.Execute("30-22=", $true, $false, $false, $false, $false, $true, 1, $false, "38-17=", 2) | Out-Null
$d.Content.Find.Execute("81+2=", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=", 2) | Out-Null
$d.Content.Find.Execute("57-30=", $true, $false, $false, $false, $false, $true, 1, $false, "8+37=", 2) | Out-Null
$d.Content.Find.Execute("88-9=", $true, $false, $false, $false, $false, $true, 1, $false, "19+30=", 2) | Out-Null
$d.Content.Find.Execute("9+45=", $true, $false, $false, $false, $false, $true, 1, $false, "81+14=", 2) | Out-Null
$d.Content.Find.Execute("9+49=", $true, $false, $false, $false, $false, $true, 1, $false, "79-67=", 2) | Out-Null
$d.Content.Find.Execute("35+29=", $true, $false, $false, $false, $false, $true, 1, $false, "61-0=", 2) | Out-Null
$d.Content.Find.Execute("77-14=", $true, $false, $false, $false, $false, $true, 1, $false, "3+16=", 2) | Out-Null
$d.Content.Find.Execute("50+41=", $true, $false, $false, $false, $false, $true, 1, $false, "35+48=", 2) | Out-Null
$d.Content.Find.Execute("12+54=", $true, $false, $false, $false, $false, $true, 1, $false, "44+53=", 2) | Out-Null
$d.Content.Find.Execute("46+43=", $true, $false, $false, $false, $false, $true, 1, $false, "8+11=", 2) | Out-Null
$d.Content.Find.Execute("32-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33-2=", 2) | Out-Null
$d.Content.Find.Execute("0+2=", $true, $false, $false, $false, $false, $true, 1, $false, "82+6=", 2) | Out-Null
